# quarterly.xlsx update: roll the quarterly data window forward by one
# quarter (drop "فصل دوم منتهی به 1399/06", add "فصل چهارم منتهی به 1401/12")
# and refresh the read_price-derived figures for the affected metric rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: quarter labels (row 8 and row 24 share the same labels) ---
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$cols = @("E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $quarters[$i]
    $ws.Range($cols[$i] + "24").Value = $quarters[$i]
}

# --- Row 16: هزینه حمل و نقل و انتقال ---
$row16 = @(-44864, 54080, 13520, 7371, 10445, 69353, 24617, 28992, 27504, 31932)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $row16[$i]
}

# --- Row 17: هزینه خدمات پس از فروش ---
$row17 = @(144953, -227113, 259312, 278674, 222401, 466912, 420387, 935435, 1222511, 3821794)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
}

# --- Row 19: هزینه مطالبات مشکوک الوصول ---
$row19 = @(2197745, 2665796, 2765959, 3359518, 3449100, 1716333, 7602949, 1771140, 4865995, 5416451)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
}

# --- Row 20: جمع ---
$row20 = @(2297834, 2492763, 3038791, 3645563, 3681946, 2252598, 8047953, 2735567, 6116010, 9270177)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
}

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت ---
$row26 = @(634, 575, 633, 575, 575, 690, 692, 692, 692, 1142)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
}

# --- Row 27: تعداد پرسنل تولیدی شرکت ---
$row27 = @(446, 504, 446, 504, 504, 517, 517, 517, 517, 2108)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
}
